# EUC_Perth_Assets.xlsx edit script
# - "4.2_Items": update LastCount/NewCount figures for several asset rows.
# - "4.2_Timestamps": fix style on the 3 existing log rows and append 10 new
#   "Add 1" log entries (rows 90-99), most recent last (descending add order
#   as per the audit log / treeview behaviour described in the commit msg).
# - "BR_Items": update LastCount/NewCount figures for a few asset rows.
# - "BR_Timestamps": append 7 new "Add 1" log entries (rows 18-24).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a literal value into a cell while (a) forcing it to be stored
# as TEXT even when it looks numeric (matches the log sheets' inlineStr
# cells) and (b) re-applying a known-good "plain" style copied from a
# neighbouring cell so the freshly-touched cell does not keep whatever
# ad-hoc number format the text coercion step added.
# ---------------------------------------------------------------------------
function Set-TextCell($cell, $value, $styleSource) {
    $cell.NumberFormat = "@"
    $cell.Value = [string]$value
    $cell.Style = $styleSource.Style
}

# ===========================================================================
# Sheet: 4.2_Items
# ===========================================================================
$items = $wb.Worksheets.Item("4.2_Items")

$items.Range("B2").Value = 3     # Desktop Mini LastCount
$items.Range("C2").Value = 6     # Desktop Mini NewCount

$items.Range("B6").Value = 51    # Laptop Charger LastCount
$items.Range("C6").Value = 53    # Laptop Charger NewCount

$items.Range("B8").Value = 0     # Monitor 24" LastCount
$items.Range("C8").Value = 2     # Monitor 24" NewCount

$items.Range("B9").Value = 0     # Monitor 34" Ultrawide LastCount
$items.Range("C9").Value = -2    # Monitor 34" Ultrawide NewCount

$items.Range("B10").Value = 4    # USB External DVD-RW Drive LastCount
$items.Range("C10").Value = 8    # USB External DVD-RW Drive NewCount

$items.Range("B11").Value = 0    # Wired Headset Poly 3325 LastCount
$items.Range("C11").Value = 4    # Wired Headset Poly 3325 NewCount

$items.Range("B14").Value = 0    # Wireless Headset Poly LastCount
$items.Range("C14").Value = 1    # Wireless Headset Poly NewCount

# ===========================================================================
# Sheet: 4.2_Timestamps
# ===========================================================================
$ts = $wb.Worksheets.Item("4.2_Timestamps")

# Rows 87-89 already hold the right values but are missing the style (s="3")
# that every other log row carries - bring them in line.
$ts.Range("A87:D89").Style = $ts.Range("A86:D86").Style

$tsStyleSrc = $ts.Range("A86")

$tsNewRows = @(
    @(90, "2023-12-11 20:35:47", "Desktop Mini",    "Add 1", "33333"),
    @(91, "2023-12-11 20:35:50", "Desktop Mini",    "Add 1", "44444"),
    @(92, "2023-12-11 20:35:54", "Desktop Mini",    "Add 1", "df333"),
    @(93, "2023-12-11 21:13:40", "Desktop Mini",    "Add 1", "12345"),
    @(94, "2023-12-11 21:17:43", "Laptop x360 G8",  "Add 1", "44444"),
    @(95, "2023-12-11 21:17:46", "Laptop x360 G8",  "Add 1", "55555"),
    @(96, "2023-12-11 21:17:49", "Laptop x360 G8",  "Add 1", "66666"),
    @(97, "2023-12-11 21:18:01", "Laptop x360 G8",  "Add 1", "44444"),
    @(98, "2023-12-11 21:18:06", "Laptop x360 G8",  "Add 1", "44444"),
    @(99, "2023-12-11 22:15:53", "Laptop x360 G8",  "Add 1", "22222")
)

foreach ($row in $tsNewRows) {
    $r = $row[0]
    Set-TextCell $ts.Range("A$r") $row[1] $tsStyleSrc
    Set-TextCell $ts.Range("B$r") $row[2] $tsStyleSrc
    Set-TextCell $ts.Range("C$r") $row[3] $tsStyleSrc
    Set-TextCell $ts.Range("D$r") $row[4] $tsStyleSrc
}

# ===========================================================================
# Sheet: BR_Items
# ===========================================================================
$britems = $wb.Worksheets.Item("BR_Items")

$britems.Range("B12").Value = 0      # Wired Keyboard LastCount
$britems.Range("C12").Value = 400    # Wired Keyboard NewCount

$britems.Range("B14").Value = 0      # Wireless Headset Poly LastCount
$britems.Range("C14").Value = -400   # Wireless Headset Poly NewCount

$britems.Range("B15").Value = 0      # Wireless Keyboard and Mouse LastCount
$britems.Range("C15").Value = 400    # Wireless Keyboard and Mouse NewCount

# ===========================================================================
# Sheet: BR_Timestamps
# ===========================================================================
$brts = $wb.Worksheets.Item("BR_Timestamps")

$brtsStyleSrc = $brts.Range("A17")

$brtsNewRows = @(
    @(18, "2023-12-11 21:18:15", "Laptop 840 G9", "Add 1", "sasdfadsdas"),
    @(19, "2023-12-11 21:18:20", "Laptop 840 G9", "Add 1", "53455534"),
    @(20, "2023-12-11 21:18:21", "Laptop 840 G9", "Add 1", "345345"),
    @(21, "2023-12-11 21:18:23", "Laptop 840 G9", "Add 1", "345345345"),
    @(22, "2023-12-11 21:18:25", "Laptop 840 G9", "Add 1", "45345345"),
    @(23, "2023-12-11 22:43:32", "Laptop 840 G9", "Add 1", "22222"),
    @(24, "2023-12-11 22:43:36", "Laptop 840 G9", "Add 1", "33333")
)

foreach ($row in $brtsNewRows) {
    $r = $row[0]
    Set-TextCell $brts.Range("A$r") $row[1] $brtsStyleSrc
    Set-TextCell $brts.Range("B$r") $row[2] $brtsStyleSrc
    Set-TextCell $brts.Range("C$r") $row[3] $brtsStyleSrc
    Set-TextCell $brts.Range("D$r") $row[4] $brtsStyleSrc
}

Write-Host "Edit complete"
